$d = $word.ActiveDocument

$replacements = @(
    @("2025-03-29 Saturday", "2025-03-30 Sunday"),
    @("397÷2=198, 1", "615÷7=87, 6"),
    @("564÷2=282, 0", "908÷8=113, 4"),
    @("514÷2=257, 0", "214÷8=26, 6"),
    @("962÷9=106, 8", "732÷9=81, 3"),
    @("518÷9=57, 5", "874÷7=124, 6"),
    @("703÷5=140, 3", "746÷8=93, 2"),
    @("715÷8=89, 3", "694÷9=77, 1"),
    @("366÷3=122, 0", "862÷2=431, 0"),
    @("333÷9=37, 0", "194÷8=24, 2"),
    @("294÷5=58, 4", "333÷9=37, 0"),
    @("470÷4=117, 2", "549÷3=183, 0"),
    @("723÷3=241, 0", "995÷5=199, 0"),
    @("128÷8=16, 0", "197÷7=28, 1"),
    @("717÷6=119, 3", "431÷2=215, 1"),
    @("639÷5=127, 4", "841÷8=105, 1"),
    @("447÷5=89, 2", "364÷8=45, 4"),
    @("898÷3=299, 1", "334÷5=66, 4"),
    @("159÷8=19, 7", "798÷5=159, 3"),
    @("960÷9=106, 6", "559÷3=186, 1"),
    @("471÷8=58, 7", "814÷3=271, 1"),
    @("444÷2=222, 0", "881÷3=293, 2"),
    @("348÷9=38, 6", "473÷8=59, 1"),
    @("497÷9=55, 2", "833÷9=92, 5"),
    @("326÷5=65, 1", "570÷4=142, 2"),
    @("729÷9=81, 0", "932÷3=310, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    [void]$range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
